$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.397.11'
$ws.Range('E2').Value = '  +4.41%  '
$ws.Range('D3').Value = '2.490.68'
$ws.Range('E3').Value = '  +2.75%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.11'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.28'
$ws.Range('E6').Value = '  +4.31%  '
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +2.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.32'
$ws.Range('E10').Value = '  +7.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0813'
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.45'
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.20'
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('D15').Value = '2.876.93'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = '2.480.75'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.846'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '47.294.08'
$ws.Range('E18').Value = '  +4.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.83'
$ws.Range('E19').Value = '  +5.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.58'
$ws.Range('E20').Value = '  +3.83%  '
$ws.Range('D21').Value = '0.0₃0939'
$ws.Range('E21').Value = '  +1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.72'
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('E23').Value = '  +6.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '251.69'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('E25').Value = '  +3.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.19'
$ws.Range('E26').Value = '  +2.00%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.26'
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.05'
$ws.Range('E29').Value = '  +4.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.24'
$ws.Range('E30').Value = '  +7.09%  '
$ws.Range('E31').Value = '  +8.76%  '
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.77'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.40'
$ws.Range('E34').Value = '  +3.71%  '
$ws.Range('E35').Value = '  +3.14%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.98'
$ws.Range('E37').Value = '  +5.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.64'
$ws.Range('E38').Value = '  +4.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.91'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('E42').Value = '  +1.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.27'
$ws.Range('E43').Value = '  +3.44%  '
$ws.Range('E44').Value = '  +2.56%  '
$ws.Range('D45').Value = '1.965.19'
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.00'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('E47').Value = '  -0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.81'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.14'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  +9.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.80'
$ws.Range('E51').Value = '  +3.55%  '
